# Added more test data to the file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the placeholder Test/IND1 sample rows with the real customer data.
# Shared-string table ends up de-duplicated & reordered by the engine based on
# write order, so we touch column C before column B on each row to reproduce
# the exact target ordering (India, John, Jose, Dennis, Raghu, Ashish appended
# after the untouched IND2..IND5 strings).
$ws.Range("C2").Value = "India"
$ws.Range("B2").Value = "John"

$ws.Range("C3").Value = "IND2"
$ws.Range("B3").Value = "Jose"

$ws.Range("C4").Value = "IND3"
$ws.Range("B4").Value = "Dennis"

$ws.Range("C5").Value = "IND4"
$ws.Range("B5").Value = "Raghu"

$ws.Range("C6").Value = "IND5"
$ws.Range("B6").Value = "Ashish"

# Column A (CustomerID) is now auto-fitted to its content.
$ws.Columns.Item(1).AutoFit()

# Final selection left on B6 (single cell) rather than the old C2:C6 range.
$ws.Range("B6").Select()
